# Weekly refresh of fruit/vegetable price data: reshuffle the daily
# records (columns D, J, K, L, M, P) among rows 2-25 according to the
# new weekly snapshot, while leaving the rest of each row intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values to copy FROM source
# row of the ORIGINAL sheet INTO destination row).
$mapping = @{
    2 = 18
    3 = 21
    4 = 8
    5 = 25
    6 = 11
    7 = 2
    8 = 6
    9 = 15
    10 = 19
    11 = 3
    12 = 9
    13 = 14
    14 = 10
    15 = 17
    16 = 5
    17 = 7
    18 = 23
    19 = 16
    20 = 13
    21 = 20
    22 = 24
    23 = 12
    24 = 4
    25 = 22
}

# Columns (by index) whose values move together as a group per row.
# D=4 (Fecha), J=10 (Volumen), K=11 (Precio minimo), L=12 (Precio maximo),
# M=13 (Precio promedio ponderado), P=16 (Precio $/Kg)
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot the current (pre-edit) values for every affected cell before
# writing anything, since several rows read from / write to each other.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2()
    }
}

# Apply the new values using the snapshot as the source of truth.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow-$c"]
    }
}
